# Apply the weekly report "no-violation" reset:
# - Bump the "Report Generated On" timestamp
# - Zero out all pricing figures (per-line Pricing column and the summary totals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report generation timestamp shown in D5
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Zero out the "Total Billed Amount" summary figure
$ws.Range("C8").Value = 0

# Zero out each line item's Pricing value (column H, rows 16-32)
for ($row = 16; $row -le 32; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}

# Zero out the grand TOTAL pricing figure
$ws.Range("H33").Value = 0
